$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "active"/"inactive" state of vendor entries across the
# All Sites / Search Engine Sites / Directory Sites / Social Sites columns.

# --- Column A (All Sites) ---
$ws.Range("A2").Value = "Bing"
$ws.Range("A3").Value = "Facebook"
$ws.Range("A4").Value = "Google"
$ws.Range("A5").Value = "Superpages"
$ws.Range("A6").Value = "TripAdvisor"
$ws.Range("A7").Value = "Yahoo"
$ws.Range("A8").Value = "Yelp"
$ws.Range("A9").Value = "YP.com"
$ws.Range("A10").ClearContents()
$ws.Range("A11").ClearContents()
$ws.Range("A12").ClearContents()
$ws.Range("A13").ClearContents()
$ws.Range("A14").ClearContents()
$ws.Range("A15").ClearContents()

# --- Column B (Search Engine Sites) ---
$ws.Range("B2").Value = "Bing"
$ws.Range("B3").Value = "Google"
$ws.Range("B5").ClearContents()

# --- Column C (Directory Sites) ---
$ws.Range("C3").Value = "Superpages"
$ws.Range("C4").Value = "Yelp"
$ws.Range("C5").Value = "YP.com"
$ws.Range("C6").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("C11").ClearContents()

# --- Selection moves to C2 ---
$ws.Range("C2").Select()
